$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Propagate D2's cell formatting (style index 4: left/center, wrap, indent)
#        down through D3:D20 so the whole "observacoes" column is consistently
#        formatted (matches the new D column cells appearing in the diff).
$ws.Range("D2").Copy()
$ws.Range("D3:D20").PasteSpecial(-4122)

# --- 2. B5: "Publicacoes/perfis..." now has its last two lines
#        ("Formulario de requisicao de voluntarios") highlighted in red,
#        since that feature (volunteers) has not shipped yet.
$fullText = "Publicações/perfis`nFormulário de inscrição em eventos`nFormulário de requisição`nde voluntários"
$ws.Range("B5").Value = $fullText
$redStart = 55
$redLength = 39
$redChars = $ws.Range("B5").Characters($redStart, $redLength)
$redChars.Font.Color = 255

# --- 3. D5: new note explaining the volunteer feature is unfinished.
$ws.Range("D5").Value = "O Voluntário ainda não foi implementado por causa de alguns BUGs em relação ao tabnavigator do mesmo."

# --- 4. Update the current selection to reflect where the user is now working.
$ws.Range("E7").Select()
